$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cosmetics")
$ws2 = $wb.Worksheets.Item("Dollar")

# Remove the barcode string from Dollar sheet row 3 (A3)
$ws2.Range("A3").ClearContents()

# Add new product row on Cosmetics sheet
$ws1.Range("A3").Value = "A1516"
$ws1.Range("B3").Value = "ADS Waterproof & Shine Lipstick"
$ws1.Range("C3").Value = 144
$ws1.Range("D3").Value = 12

$ws1.Range("A1:D3").SetPhonetic()

$ws1.Activate()
$ws1.Range("C29").Select()
